# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at the top of the data (row 7),
# pushing the existing rows 7-24 down to 8-25, and populate the new
# row 7 with the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = 45250
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100102
$ws.Range("H7").Value = "Cítricos"
$ws.Range("I7").Value = 100102006
$ws.Range("J7").Value = "Pomelo"
$ws.Range("K7").Value = "Start Ruby"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 150
$ws.Range("N7").Value = 17000
$ws.Range("O7").Value = 17000
$ws.Range("P7").Value = 17000
$ws.Range("Q7").Value = "$/caja 14 kilos empedrada"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1214
$ws.Range("T7").Value = 14
